# PIList.xlsx edit: add "Old New Faces" / "Math bio" PI entries and
# re-sort the table by Last Name, First Name (mirrors the sheet's sortState).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing data occupies rows 2..45 (44 PI rows under the header row 1).
$lastExistingRow = 45

# New PI rows to add (Last Name, First Name, Feature).
$newRows = @(
    @("Ding", "Fanguayn", "Old New Faces"),
    @("Komarova", "Natalia", "Math bio"),
    @("Kvon", "Evgeny", "Old New Faces"),
    @("Lee", "Grace", "Old New Faces"),
    @("Li", "Wei", "Old New Faces"),
    @("Rodriguez-Verdugo", "Alejandra", "Old New Faces"),
    @("Seldin", "Marcu", "Old New Faces"),
    @("Shi", "Xiaoyu", "Old New Faces"),
    @("Siryaporn", "Albert", "Math bio"),
    @("Wodarz", "Dominik", "Math bio"),
    @("Xu", "Xiangmin", "Old New Faces"),
    @("Yu", "Jin", "Old New Faces")
)

$row = $lastExistingRow + 1
foreach ($entry in $newRows) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}

$lastRow = $row - 1

# Step 1: perform the actual re-sort of the data range (A2:C<lastRow>) using
# Last Name then First Name as keys (ascending) so ties between existing and
# newly-added same-surname rows land in the correct order.
$sortRange = $ws.Range("A2:C" + $lastRow)
$key1 = $ws.Range("A2:A" + $lastRow)
$key2 = $ws.Range("B2:B" + $lastRow)
$sortRange.Sort($key1, 1, $key2, $null, 1)

# Step 2: refresh the worksheet's remembered sort state (single condition on
# column A, like the original workbook) so it points at the new extended
# range instead of the old one. The original file's sortState range extended
# two rows past the actual data (A2:C47 vs dimension C45); the edited file's
# sortState similarly extends one row past the new data (A2:C58 vs C57), so
# we reproduce that same remembered range here.
$sortStateLastRow = $lastRow + 1
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A2:A" + $lastRow))
$sortObj.SetRange($ws.Range("A2:C" + $sortStateLastRow))
$sortObj.Header = 0
$sortObj.Apply()

# Update the active cell selection to match the target workbook state.
$ws.Range("E12").Select()

Write-Host ("Last row after sort: " + $lastRow)
